# Word COM-interop script implementing the tracked changes.
#
# Summary of the source diff:
#   1) "Contents" (TOC heading) gets wrapped in spellcheck proofErr marks -
#      no visible text changes (Word-internal spell-check bookkeeping only,
#      not reachable through the Word object model - skipped).
#   2) The "... Suralovi, vinari z obce Petrov. D" run is re-split into more
#      runs (and wrapped in proofErr marks) but the concatenated text is
#      byte-for-byte identical before/after - nothing to type.
#   3) The validity paragraph gains a new clause: "... porovnavaly se
#      ziskanymi udaji" becomes "... porovnavaly se ziskanymi udaji za
#      poslednich 8 let. Tyto poznatky a dalsi navrhnute experimenty byly
#      nasledne konzultovany s Petrem Suralem." - this is the one real
#      content edit, applied below.
#   4) A cached TOC page-number field result flips from "2" to "3" - this
#      is a calculated field value normally refreshed by Word's pagination
#      engine (F9 / update-fields-on-print), not something an editor sets
#      by typing.
#
$d = $word.ActiveDocument

# --- Paragraph "Overovani validity modelu probihalo prubezne ..." ---------
# Locate it by scanning paragraphs for the distinctive lead-in text instead
# of a hard-coded index, so the script is robust to minor structural drift.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $s = $para.Range.Start
    $e = $para.Range.End
    if ($e -le $s) { continue }
    # Trim the trailing paragraph-mark position: re-reading through a fresh
    # Range keyed on (Start, End-1) avoids the "whole paragraph re-insert"
    # quirk some hosts exhibit when a Range spans its own end-of-paragraph.
    $probe = $d.Range($s, $e - 1)
    $t = $probe.Text
    if ($t -and $t.Contains("porovnávaly se získanými údaji")) {
        $target = $probe
        break
    }
}

if ($target -ne $null) {
    $old = $target.Text
    $oldTail = " a následně byly konzultovány s Petrem Šuralem."
    $newTail = " za posledních 8 let. Tyto poznatky a další navrhnuté experimenty byly následně konzultovány s Petrem Šuralem."
    if ($old.Contains($oldTail)) {
        $target.Text = $old.Replace($oldTail, $newTail)
    }
}
